# Import of settings implemented
#
# Wires the add-in's "Settings" sheet up to an already-imported
# configuration: external-source usage now pulls straight from the
# TitleBlockData sheet (starting at B2), the settings-sheet reference
# points at the real "Settings" tab, the various boolean-looking flags
# are normalised to upper-case TRUE text (matching the B2/B6/.../B19
# dropdown's literal "True,False" list values) with the import-related
# switches turned on, and the stray MapNoSheets mapping is cleared out.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TitleBlockData")
$ws2 = $wb.Worksheets.Item("Settings")

# Helper: typing the literal word TRUE/FALSE straight into a cell makes
# Excel auto-convert it to a real Boolean. These settings are read back
# as plain text (they feed a "True,False" validation dropdown), so build
# the word via a formula in a scratch cell and paste-special *values
# only* into each target cell - exactly like pasting in text from
# another source keeps it text instead of re-parsing it as a literal.
function Set-TrueText([string]$addr) {
    $ws2.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

$scratch = $ws2.Range("D1")
$scratch.Formula = "=""TR""&""UE"""
$scratch.Copy()

# --- Row 2: "Use external source" -> "UseExternalSource" / True -> TRUE
$ws2.Range("A2").Value = "UseExternalSource"
Set-TrueText "B2"

# --- Row 4: SheetName -> point at the real data sheet
$ws2.Range("B4").Value = "TitleBlockData"

# --- Row 5: StartCell -> B2
$ws2.Range("B5").Value = "B2"

# --- Row 6: AutoFillTitleBlock -> TRUE
Set-TrueText "B6"

# --- Row 7: ImportSettingsXL -> TRUE
Set-TrueText "B7"

# --- Row 8: SheetName_Settings -> the actual Settings sheet/tab name
$ws2.Range("B8").Value = "Settings"

# --- Row 10: UseFileName -> TRUE
Set-TrueText "B10"

# --- Row 11: DrwNrFieldName, drop the leading apostrophe from the value
$ws2.Range("B11").Value = "DRAWING_TITLE"

# --- Row 15: MapNoSheets -> cleared
$ws2.Range("B15").Value = ""

# --- Rows 16-19: Include* flags -> TRUE
Set-TrueText "B16"
Set-TrueText "B17"
Set-TrueText "B18"
Set-TrueText "B19"

# Clean up the scratch helper cell used to build the TRUE text.
$scratch.Clear()

# Settings sheet column A width 24 -> 23 characters.
# (ColumnWidth's setter adds Excel's usual ~0.83-character padding vs.
# the raw stored width, so dial the input back to land on exactly 23.)
$ws2.Columns.Item(1).ColumnWidth = 22.1666666666667

# Make "Settings" the active sheet/tab (activeTab 0 -> 1).
$ws2.Activate()
